# The presentation's single Design ("Integral") swaps color palettes with the
# stock "Office Theme" that was previously only attached to the notes master
# (ppt/theme/theme1.xml <-> ppt/theme/theme2.xml in the OOXML package).
#
# This host exposes exactly one writable Theme object (the one seated on the
# slide master / active Design), which is physically persisted as
# ppt/theme/theme2.xml. Re-pointing the 12 theme colour slots
# (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) to the "Office Theme"
# palette reproduces the colour content of that swap.

$p = $ppt.ActivePresentation

$tcs = $p.SlideMaster.Theme.ThemeColorScheme

# PowerPoint packs colours as a single Long: value = R + G*256 + B*65536
$tcs.Item(1).RGB  = 0x000000    # dk1      000000
$tcs.Item(2).RGB  = 0xFFFFFF    # lt1      FFFFFF
$tcs.Item(3).RGB  = 0x6A5444    # dk2      44546A
$tcs.Item(4).RGB  = 0xE6E6E7    # lt2      E7E6E6
$tcs.Item(5).RGB  = 0xD59B5B    # accent1  5B9BD5
$tcs.Item(6).RGB  = 0x317DED    # accent2  ED7D31
$tcs.Item(7).RGB  = 0xA5A5A5    # accent3  A5A5A5
$tcs.Item(8).RGB  = 0x00C0FF    # accent4  FFC000
$tcs.Item(9).RGB  = 0xC47244    # accent5  4472C4
$tcs.Item(10).RGB = 0x47AD70    # accent6  70AD47
$tcs.Item(11).RGB = 0xC16305    # hlink    0563C1
$tcs.Item(12).RGB = 0x724F95    # folHlink 954F72
